$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SQL queries in column B (rows 2-7) and C2 all joined on the generic
# ".id" columns (std.id / prt.id). The schema was updated to use the
# explicit natural keys (study_id / participant_id), so every query needs
# its JOIN predicates rewritten accordingly.
$cellsToFix = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cellsToFix) {
    $text = $ws.Range($addr).Value2
    if ($text -eq $null) { continue }

    $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    $ws.Range($addr).Value = $text
}

# Widen column C (it's no longer an exact "best fit" width - the user
# resized it manually), and move the active selection/viewport up to B2.
$ws.Columns.Item(3).ColumnWidth = 68.3
$ws.Range("B2").Select()
